# 20.03.06 PM 20:24 Lot 속성 부여 로직 작성 중
# Swap the B/C values of the first two rows in each 5-row grade block
# (rows 11-50), then clear out the duplicated block (rows 51-90) and
# move the viewport/selection to reflect where editing left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(11, 16, 21, 26, 31, 36, 41, 46)
foreach ($r1 in $rowPairs) {
    $r2 = $r1 + 1

    $bVal1 = $ws.Cells.Item($r1, 2).Value2
    $cVal1 = $ws.Cells.Item($r1, 3).Value2
    $bVal2 = $ws.Cells.Item($r2, 2).Value2
    $cVal2 = $ws.Cells.Item($r2, 3).Value2

    $ws.Cells.Item($r1, 2).Value2 = $bVal2
    $ws.Cells.Item($r1, 3).Value2 = $cVal2
    $ws.Cells.Item($r2, 2).Value2 = $bVal1
    $ws.Cells.Item($r2, 3).Value2 = $cVal1
}

# Remove the duplicated tail block (rows 51-90), shrinking the used range
# back down to A1:C50.
$ws.Range("A51:C90").ClearContents()

# Reflect the scrolled viewport / active selection at the end of the edit.
$ws.Application.ActiveWindow.ScrollRow = 35
$ws.Range("K41").Select()
